$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.193862795829773
$ws.Range("B1").Value = 3.853585481643677
$ws.Range("C1").Value = 3.154757499694824
$ws.Range("D1").Value = 2.529603481292725
$ws.Range("E1").Value = 1.365979909896851
